# Remove comments levenberg_marquadt: insert a new "Weights_raw" column
# (raw/unconverted final weights) before the existing "Final_Weights"
# column, and refresh a few Initial_Weights / Final_Weights values that
# changed as a result of no longer rounding/formatting those arrays.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D; everything from D onward shifts right by one.
$ws.Columns("D:D").Insert()

# New header for the inserted column.
$ws.Range("D1").Value = "Weights_raw"

# New "Weights_raw" column values (previously not present in the sheet).
$ws.Range("D2").Value = "[ 5.0081046   5.73419945 10.59963756]"
$ws.Range("D3").Value = "[ 5.0081046   5.73419945 10.59963756]"
$ws.Range("D4").Value = "[371.16064774  47.70853941  13.11971189]"
$ws.Range("D5").Value = "[371.16064774  47.70853941  13.11971189]"
$ws.Range("D6").Value = "[288.11203268  48.82467957  36.99057293]"
$ws.Range("D7").Value = "[288.11203268  48.82467957  36.99057293]"

# Updated Initial_Weights values (column C) for rows 3, 5, 7.
$ws.Range("C3").Value = "[0.04271399961972451, 0.008952370098035345, 0.07679250311927917]"
$ws.Range("C5").Value = "[0.04271399961972451, 0.008952370098035345, 0.07679250311927917]"
$ws.Range("C7").Value = "[0.04271399961972451, 0.008952370098035345, 0.07679250311927917]"

# Updated Final_Weights values (now column E, was D) for rows 3, 5, 7.
$ws.Range("E3").Value = "[5.008104595025083, 5.73419944808188, 10.599637559381778]"
$ws.Range("E5").Value = "[5.008104595001941, 5.734199448076006, 638.9034182886581]"
$ws.Range("E7").Value = "[5.008104595020175, 5.734199448068985, 10.599637559494475]"
